$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text updates (values that Excel will not auto-convert to numbers)
$plainUpdates = @{
    "D2" = "71.583.04"
    "E2" = "  +1.51%  "
    "D3" = "3.830.28"
    "E4" = "  +0.05%  "
    "E5" = "  +4.88%  "
    "E6" = "  +3.24%  "
    "D7" = "3.827.95"
    "E7" = "  +0.45%  "
    "E8" = "  -0.03%  "
    "E9" = "  +0.11%  "
    "E10" = "  +1.41%  "
    "E11" = "  +4.98%  "
    "E12" = "  -0.48%  "
    "E13" = "  +5.97%  "
    "E14" = "  +2.15%  "
    "D15" = "4.469.79"
    "E15" = "  +0.44%  "
    "D16" = "3.821.62"
    "E16" = "  +0.39%  "
    "D17" = "71.495.43"
    "E17" = "  +1.52%  "
    "E18" = "  +0.67%  "
    "E19" = "  +1.22%  "
    "E20" = "  +0.46%  "
    "E21" = "  +2.36%  "
    "E22" = "  +2.75%  "
    "E23" = "  +0.63%  "
    "E24" = "  +2.17%  "
    "E25" = "  +0.08%  "
    "E26" = "  +1.17%  "
    "E27" = "  +1.99%  "
    "E28" = "  +1.66%  "
    "D29" = "3.974.26"
    "E29" = "  +0.26%  "
    "E30" = "  +0.14%  "
    "E31" = "  +9.99%  "
    "E32" = "  -0.02%  "
    "E33" = "  +2.99%  "
    "E34" = "  +0.34%  "
    "E35" = "  +1.82%  "
    "E36" = "  +1.79%  "
    "D37" = "3.777.68"
    "E38" = "  -0.03%  "
    "E39" = "  +1.58%  "
    "E40" = "  +15.26%  "
    "E41" = "  -0.40%  "
    "E42" = "  +1.41%  "
    "E43" = "  +3.23%  "
    "E44" = "  -0.05%  "
    "E46" = "  +8.44%  "
    "E47" = "  +3.86%  "
    "E48" = "  -2.36%  "
    "E49" = "  +1.39%  "
    "E50" = "  +0.64%  "
    "E51" = "  +2.32%  "
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# Numeric-looking text updates (must be forced to stay text, matching the
# original inline-string cell type instead of being parsed as a number)
$textForcedUpdates = @{
    "D4" = "0.999"
    "D5" = "701.06"
    "D6" = "174.53"
    "D9" = "0.527"
    "D11" = "7.43"
    "D14" = "36.67"
    "D18" = "17.83"
    "D19" = "7.26"
    "D22" = "488.15"
    "D24" = "84.73"
    "D25" = "0.0000145"
    "D26" = "12.40"
    "D27" = "10.56"
    "D28" = "2.16"
    "D33" = "7.61"
    "D34" = "29.77"
    "D35" = "0.182"
    "D36" = "9.34"
    "D41" = "3.44"
    "D44" = "0.999"
    "D46" = "0.000313"
    "D47" = "162.93"
    "D48" = "44.75"
    "D49" = "48.76"
    "D51" = "8.71"
}
foreach ($ref in $textForcedUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textForcedUpdates[$ref]
    $cell.ClearFormats()
}
